$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.435.39'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.062.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.43%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.93'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.25'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +6.31%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.59'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0760'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.368.03'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.33'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.03'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.775'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.16'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.044.48'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '37.360.69'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.20'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +15.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.19'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.29%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '225.71'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.28%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.43'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.40'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.60'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.91%  '
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.93'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.87%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.46'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.129'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.12'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.117'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.58'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +6.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0617'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.55'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.27%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.92'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.43%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.30'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.69'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +13.31%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.34%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.467.24'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.30'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +6.91%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0939'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.17'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.43%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0211'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.60'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.02'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.22'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.24%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.06%  '
